$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the O column formula: remove the J (fuhep) factor from the calculation.
# New formula: =K{row}/120/21.4*1000
$ws.Range("O2:O119").Formula = "=K2/120/21.4*1000"

# Apply a numeric format (0.00) to the O column formula cells to match the new style.
$ws.Range("O2:O119").NumberFormat = "0.00"

# Update sheet view: remove the saved topLeftCell / old selection, select O2:O119 with active cell O2.
$ws.Range("O2:O119").Select()
